$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B. This shifts the existing
# "Jun_13" column (B) to D, and the existing "Jun_10" column (C) to E -
# replicating the rolling-window report picking up two newer report
# dates (Jun_17, Jun_15) while keeping the older ones around.
$ws.Columns("B:C").Insert()

# The two freshly inserted columns (B and C) default every data row to
# "UN" (unrated), matching how every other date column starts out.
$ws.Range("B2:C27").Value = "UN"

# Header row: newest date goes in B, next newest in C. D1/E1 already
# hold the shifted "Jun_13"/"Jun_10" headers from the Insert above.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Match the original column C formatting (custom width 8.0) across all
# three of the now-adjacent columns C, D and E.
$ws.Columns("C").ColumnWidth = 7.166666666666667
$ws.Columns("D").ColumnWidth = 7.166666666666667
$ws.Columns("E").ColumnWidth = 7.166666666666667
